# Commit: Committing index.html for Game {63960141}, Stage 2
# Moves the 4 existing winner-announcement rows from rows 7-10 down to
# rows 11-14 (new timestamps), leaves rows 7-10 blank again, updates the
# 4 "winning group" texts, and appends 4 new blank rows (111-114) at the
# bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Move the data currently sitting in rows 7-10 (timestamp + group
#    text, with their number/time formatting) down to rows 11-14.
#    Copy also carries over the per-cell styles (s="2" / s="3").
# ---------------------------------------------------------------------
$ws.Range("A7:B10").Copy($ws.Range("A11:B14"))

# ---------------------------------------------------------------------
# 2) Blank out the old rows 7-10 again: remove their values, and reset
#    their formatting back to the plain "empty data row" style (the same
#    style already used by columns C:H on those rows, i.e. style index 1).
# ---------------------------------------------------------------------
$ws.Range("A7:B10").ClearContents()
$ws.Range("C7:C10").Copy()
$ws.Range("A7:B10").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 3) Update the timestamps for the (now relocated) submissions.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 45285.022210648145
$ws.Range("A12").Value = 45285.022337962961
$ws.Range("A13").Value = 45285.022418981483
$ws.Range("A14").Value = 45285.022499999999

# ---------------------------------------------------------------------
# 4) Update the winning-group text for each submission.
# ---------------------------------------------------------------------
$ws.Range("B11").Value = "הקבוצה של: אור, שרי"
$ws.Range("B12").Value = "הקבוצה של: המפקד, הקשבי"
$ws.Range("B13").Value = "הקבוצה של: איי, עמרי"
$ws.Range("B14").Value = "הקבוצה של: דור, גור"

# ---------------------------------------------------------------------
# 5) Append 4 new blank rows (111-114) after the current last row (110),
#    matching the same style/formatting and row height as the existing
#    blank rows.
# ---------------------------------------------------------------------
$ws.Range("A110:H110").Copy($ws.Range("A111:H111"))
$ws.Range("A110:H110").Copy($ws.Range("A112:H112"))
$ws.Range("A110:H110").Copy($ws.Range("A113:H113"))
$ws.Range("A110:H110").Copy($ws.Range("A114:H114"))
$ws.Range("A111:H114").ClearContents()
$ws.Rows("111:114").RowHeight = 15.75

# ---------------------------------------------------------------------
# 6) Move the active-cell selection to D9 (matches the author's saved
#    cursor position after editing).
# ---------------------------------------------------------------------
$ws.Range("D9").Select() | Out-Null
